# Trade #6 closed at 2026-02-17 13:08:04 - unknown UNKNOWN +0.000%
#
# Updates the live trading results workbook after a new (6th) trade on the
# MarketMaking strategy was closed:
#   - Summary sheet totals refreshed
#   - Strategy Status row for MarketMaking refreshed
#   - New trade row appended to both "All Trades" and "MarketMaking" sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("B3").Value = 1199.59   # Current Capital
$summary.Range("B4").Value = -0.41     # Total P&L $
$summary.Range("B5").Value = -1.37     # Total P&L %
$summary.Range("B6").Value = 6         # Total Trades
$summary.Range("B7").Value = 2         # Winning Trades
$summary.Range("B9").Value = 33.33     # Win Rate %

# ---------------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")

$status.Range("C4").Value = 99.59      # Capital
$status.Range("D4").Value = 6          # Trades
$status.Range("E4").Value = -0.41      # P&L $
$status.Range("F4").Value = -0.41      # P&L %
$status.Range("G4").Value = 33.33      # Win Rate %

# ---------------------------------------------------------------------------
# 3. Append new trade row (row 7) to "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------------------
function Add-TradeRow($ws) {
    $ws.Cells.Item(7, 1).Value = 6
    $ws.Cells.Item(7, 2).NumberFormat = "@"
    $ws.Cells.Item(7, 2).Value = "2026-02-17"
    $ws.Cells.Item(7, 3).Value = "13:07:58"
    $ws.Cells.Item(7, 4).Value = "MarketMaking"
    $ws.Cells.Item(7, 5).Value = "UP"
    $ws.Cells.Item(7, 6).Value = 0.14
    $ws.Cells.Item(7, 7).Value = 0.21
    $ws.Cells.Item(7, 8).Value = "CLOSED"
    $ws.Cells.Item(7, 9).Value = 50
    $ws.Cells.Item(7, 10).Value = 0.07000000000000001
    $ws.Cells.Item(7, 11).Value = 99.59
    $ws.Cells.Item(7, 12).Value = 0
    $ws.Cells.Item(7, 13).Value = 0
    $ws.Cells.Item(7, 14).Value = 0.6
    $ws.Cells.Item(7, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(7, 16).Value = "early_exit"
    $ws.Cells.Item(7, 17).Value = 0.13
}

Add-TradeRow($wb.Worksheets.Item("All Trades"))
Add-TradeRow($wb.Worksheets.Item("MarketMaking"))
